# Apply updated cryptocurrency market data (price & volume change) per source refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.145.41'
$ws.Range("E2").Value = '  -0.58%  '

$ws.Range("D3").Value = '1.778.01'
$ws.Range("E3").Value = '  -2.55%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.12'
$ws.Range("E5").Value = '  -1.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.551'
$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.09'
$ws.Range("E8").Value = '  +1.39%  '

$ws.Range("E9").Value = '  -1.49%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0657'
$ws.Range("E10").Value = '  -2.53%  '

$ws.Range("E11").Value = '  +0.02%  '

$ws.Range("D12").Value = '2.032.04'
$ws.Range("E12").Value = '  -2.58%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.17'
$ws.Range("E13").Value = '  +7.21%  '

$ws.Range("D14").Value = '1.776.55'
$ws.Range("E14").Value = '  -2.85%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.626'
$ws.Range("E15").Value = '  -3.33%  '

$ws.Range("D16").Value = '34.117.93'
$ws.Range("E16").Value = '  -0.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.21'
$ws.Range("E17").Value = '  -1.89%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.69'
$ws.Range("E18").Value = '  -1.61%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '255.03'
$ws.Range("E19").Value = '  -1.45%  '

$ws.Range("D20").Value = '0.0₃0739'
$ws.Range("E20").Value = '  -2.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.37'
$ws.Range("E22").Value = '  -2.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.20'
$ws.Range("E23").Value = '  -3.40%  '

$ws.Range("E24").Value = '  -3.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.62'
$ws.Range("E25").Value = '  -1.90%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.39'
$ws.Range("E26").Value = '  -1.52%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.00'
$ws.Range("E27").Value = '  -2.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.114'
$ws.Range("E28").Value = '  -1.42%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.77'
$ws.Range("E30").Value = '  -3.53%  '

$ws.Range("E31").Value = '  -1.77%  '

$ws.Range("E32").Value = '  -1.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.59'
$ws.Range("E33").Value = '  +0.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.85'
$ws.Range("E34").Value = '  +2.88%  '

$ws.Range("D35").Value = '1.436.24'
$ws.Range("E35").Value = '  -7.27%  '

$ws.Range("E36").Value = '  -3.83%  '

$ws.Range("E37").Value = '  -1.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.624'
$ws.Range("E38").Value = '  -1.42%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.86'
$ws.Range("E39").Value = '  +1.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '82.77'
$ws.Range("E40").Value = '  -2.35%  '

$ws.Range("E41").Value = '  +0.51%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.886'
$ws.Range("E42").Value = '  -3.67%  '

$ws.Range("E43").Value = '  -5.29%  '

$ws.Range("E44").Value = '  -2.66%  '

$ws.Range("E45").Value = '  -2.00%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.82'
$ws.Range("E46").Value = '  +0.80%  '

$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '1.933.93'
$ws.Range("E47").Value = '  -2.81%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.21'
$ws.Range("E48").Value = '  -0.65%  '

$ws.Range("E49").Value = '  +0.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '98.20'
$ws.Range("E50").Value = '  +0.73%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '49.54'
$ws.Range("E51").Value = '  -6.21%  '
